$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) A few pre-existing cells pick up a highlight (style) that didn't
#    exist before: reuse the "red" fill already used elsewhere
#    (style index 6) by copying *formats only* from a cell that
#    already carries it, so no new style/fill gets minted.
# ------------------------------------------------------------------
$ws.Range("D2").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$ws.Range("C17:E17").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# ------------------------------------------------------------------
# 2) Append a new ROLE/ROUTE table for the "SCHOOL" route (rows 21-24),
#    mirroring the layout of the existing blocks above it.
#    Clone formatting from the previous block (rows 16-19) first so
#    every cell lands on the exact same style indices, then overwrite
#    the text values on top (value writes do not disturb style here).
# ------------------------------------------------------------------
$ws.Range("A16:J19").Copy() | Out-Null
$ws.Range("A21:J24").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# Merge the ROUTE label cell across the new block, same as every other
# block. Merge the (still blank-of-that-style) cells first, *then* stamp
# the final format on top -- merging already-styled cells makes Excel
# split the border into top/middle/bottom variants, which the source
# blocks above don't have.
$ws.Range("B22:B24").Merge() | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B22:B24").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# Row 22/23 (STUDENT/TEACHER) mirror the highlighted CREATE/DELETE/UPDATE
# cells used elsewhere (e.g. row 7) rather than the plain copied style.
$ws.Range("C7:E7").Copy() | Out-Null
$ws.Range("C22:E22").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$ws.Range("C22:E22").Copy() | Out-Null
$ws.Range("C23:E23").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# Now fill in the text values.
$ws.Range("A21").Value = "ROLE"
$ws.Range("B21").Value = "ROUTE"
$ws.Range("C21").Value = "CREATE 1"
$ws.Range("D21").Value = "DELETE 1"
$ws.Range("E21").Value = "UPDATE 1"
$ws.Range("F21").Value = "READ 1"
$ws.Range("G21").Value = "BULK CREATE"
$ws.Range("H21").Value = "BULK DELETE"
$ws.Range("I21").Value = "BULK UPDATE"
$ws.Range("J21").Value = "BULK READ"

$ws.Range("A22").Value = "STUDENT"
$ws.Range("B22").Value = "SCHOOL"
$ws.Range("A23").Value = "TEACHER"
$ws.Range("A24").Value = "ADMIN"

# Reset selection back to the top-left (clears the stale H24 selection).
$ws.Range("A1").Select() | Out-Null
